# The "Warmup Plan" sheet had two stray leftover label rows at the very
# top (A1="Properties", A2="Value") that don't belong to the actual
# table (whose real header - Phase/Run/Gmail/Yahoo/Microsoft/Adobe/
# Roadrunner/Others - starts on row 3). Remove those two rows so the
# table header becomes row 1 and the data shifts up accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warmup Plan")

$ws.Rows("1:2").Delete()

# Restore the selection to the (now) header row, matching how the sheet
# was left selected after the cleanup.
$ws.Rows("1:1").Select()
